# Applies the "Added 2030 DH units" + FI00 DH demand correction-factor edit
# to the demand.xlsx workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Correction factor: FI00 / hydrogen / National Trends / 2025 demand
#    drops from 400 to 300 (row 2, column E).
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 5).Value = 300

# ---------------------------------------------------------------------
# 2. Add the new 2030 DH unit rows (21-26) right below the existing data.
#    Node, Inputnodestub, Scenario, Year, Value
# ---------------------------------------------------------------------
$newRows = @(
    @("FI00", "hydrogen", "Distributed Energy", 2030, 500),
    @("SE01", "hydrogen", "Distributed Energy", 2030, 500),
    @("SE02", "hydrogen", "Distributed Energy", 2030, 250),
    @("NON1", "hydrogen", "Distributed Energy", 2030, 100),
    @("DE00", "hydrogen", "Distributed Energy", 2030, 2000),
    @("FR00", "hydrogen", "Distributed Energy", 2030, 1000)
)

$r = 21
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Re-apply the AutoFilter over the extended range (A1:E20), this time
#    filtering on "National Trends" instead of "Distributed Energy".
#    Removing the previous filter first lets the new range stick.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:E20").AutoFilter(3, @("National Trends"))

# Keep the hidden "_FilterDatabase" name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "demand!_FilterDatabase") {
        $n.RefersTo = '=demand!$A$1:$E$20'
    }
}

# ---------------------------------------------------------------------
# 4. Match the selection left behind by the editor.
# ---------------------------------------------------------------------
$ws.Range("E37").Select()
